# Apply the "Modified Work Order templates and other fixes" edit:
#  1. Delete the unused "Sheet1" worksheet (keep only "Create WO").
#  2. Update the lot-track project value on "Create WO"!B2.
#  3. Move the active selection from A3:XFD3 (row 3) to E2.

$wb = $excel.ActiveWorkbook

# Avoid any "are you sure you want to delete" prompt when removing the sheet.
$excel.DisplayAlerts = $false

$wsCreateWO = $wb.Worksheets.Item("Create WO")
$wsSheet1   = $wb.Worksheets.Item("Sheet1")

$wsSheet1.Delete()

# Update the project/lot-track description text in row 2.
$wsCreateWO.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Reset the selection on the remaining sheet to cell E2.
$wsCreateWO.Activate()
$wsCreateWO.Range("E2").Select()
